$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of (cellAddress, newValue) pairs derived from the diff.
$edits = @(
    @("D2", "65.800.91"),
    @("E2", "  -2.28%  "),
    @("D3", "3.752.80"),
    @("E3", "  +1.04%  "),
    @("E4", "  -0.08%  "),
    @("D5", "404.04"),
    @("E5", "  -4.19%  "),
    @("D6", "132.17"),
    @("E6", "  +0.24%  "),
    @("D7", "3.743.07"),
    @("E7", "  +1.05%  "),
    @("E8", "  -6.16%  "),
    @("E9", "  +0.05%  "),
    @("D10", "0.718"),
    @("E10", "  -6.57%  "),
    @("D11", "0.166"),
    @("E11", "  -9.95%  "),
    @("D12", "0.0000355"),
    @("E12", "  -11.85%  "),
    @("D13", "40.68"),
    @("E13", "  -6.02%  "),
    @("D14", "4.360.14"),
    @("E14", "  +1.29%  "),
    @("D15", "9.74"),
    @("E15", "  -5.60%  "),
    @("D16", "14.77"),
    @("E16", "  +12.80%  "),
    @("E17", "  -1.53%  "),
    @("D18", "3.764.29"),
    @("E18", "  +1.52%  "),
    @("D19", "19.43"),
    @("E19", "  -6.79%  "),
    @("D20", "65.988.94"),
    @("E20", "  -2.05%  "),
    @("E21", "  -6.43%  "),
    @("D22", "408.72"),
    @("E22", "  -9.59%  "),
    @("D23", "14.30"),
    @("E23", "  -9.45%  "),
    @("D24", "84.91"),
    @("E24", "  -5.47%  "),
    @("D25", "3.04"),
    @("E25", "  -4.49%  "),
    @("D26", "5.69"),
    @("E26", "  +14.20%  "),
    @("D27", "35.78"),
    @("E27", "  -5.98%  "),
    @("E28", "  -7.19%  "),
    @("D29", "9.32"),
    @("E29", "  -8.64%  "),
    @("D30", "12.34"),
    @("E30", "  -2.50%  "),
    @("B31", "Hedera"),
    @("C31", "https://coinranking.com/coin/jad286TjB+hedera-hbar"),
    @("D31", "0.119"),
    @("E31", "  -3.16%  "),
    @("B32", "Toncoin"),
    @("C32", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"),
    @("D32", "2.67"),
    @("E32", "  -2.38%  "),
    @("D33", "7.35"),
    @("E33", "  -0.78%  "),
    @("E34", "  -6.73%  "),
    @("D35", "38.93"),
    @("E35", "  -7.49%  "),
    @("D36", "1.00"),
    @("E36", "  +0.06%  "),
    @("D37", "54.88"),
    @("E37", "  -2.95%  "),
    @("D38", "0.0₃0730"),
    @("E38", "  -6.08%  "),
    @("E39", "  -7.91%  "),
    @("E40", "  -9.67%  "),
    @("D41", "0.998"),
    @("E41", "  -0.08%  "),
    @("E42", "  -8.94%  "),
    @("D43", "27.24"),
    @("E43", "  -1.82%  "),
    @("B44", "ApeXProtocol"),
    @("C44", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"),
    @("D44", "3.18"),
    @("E44", "  +20.00%  "),
    @("B45", "Monero"),
    @("C45", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("D45", "144.66"),
    @("E45", "  -2.39%  "),
    @("D46", "3.24"),
    @("E46", "  -5.45%  "),
    @("E47", "  -2.28%  "),
    @("D48", "2.61"),
    @("E48", "  -2.41%  "),
    @("D49", "4.25"),
    @("E49", "  -4.29%  "),
    @("E50", "  -5.29%  "),
    @("D51", "0.293"),
    @("E51", "  -5.26%  ")
)

# Column D frequently contains values that look like plain numbers (e.g. "9.74"),
# but in this workbook they must remain text (inline/shared strings), matching
# the original file's cell typing. Temporarily force column D to a Text number
# format while assigning the values, then clear the temporary formatting again
# so the workbook's style table stays as close as possible to its original state.
$colDRange = $ws.Range("D2:D51")
$originalFormat = $colDRange.NumberFormat
$colDRange.NumberFormat = "@"

foreach ($edit in $edits) {
    $addr = $edit[0]
    $value = $edit[1]
    $ws.Range($addr).Value = $value
}

$colDRange.ClearFormats()
